$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at BE. This shifts the existing "nom" (BE) and
# "url_produit" (BF) columns one position to the right (BF and BG),
# and copies formatting from the previous column into the new BE column.
$ws.Range("BE1").EntireColumn.Insert()

# New header cell: latest scrape timestamp for this price-history column.
$ws.Range("BE1").Value2 = "2026-01-30 07:33:59"

# Seed the new price column with the most recent known price, copied
# from column BD, for every data row that currently holds a price (rows 2-80).
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $bdCell = $ws.Cells.Item($r, 56)   # column BD = 56
    $priceVal = $bdCell.Value2
    if ($priceVal -ne $null -and $priceVal -ne "") {
        $ws.Cells.Item($r, 57).Value2 = $priceVal   # column BE = 57
    }
}
